# Rename the "AddressBook" class box to "DeskBoard" on the Logic Component
# Class Diagram slide (the diagram's top line reads "AddressBook" / "Parser"
# as two separate paragraphs inside the same shape - only the first
# paragraph's text changes).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tf = $shp.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            if ($tr.Text -eq "AddressBook`rParser") {
                $tr.Paragraphs(1).Runs(1).Text = "DeskBoard"
            }
        }
    }
}
